$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# For D-column values that look like plain numbers (e.g. "311.31"), the
# cell's number format is temporarily switched to Text so Excel keeps the
# exact original text formatting (trailing zeros, etc.) instead of
# auto-converting the entry to a numeric value. The style is restored to
# Normal afterwards so no visible formatting changes are introduced.
$ws.Range("D2").Value = '45.045.26'
$ws.Range("E2").Value = '  +2.12%  '
$ws.Range("D3").Value = '2.357.01'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.31'
$ws.Range("E5").Value = '  -1.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.76'
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.606'
$ws.Range("E9").Value = '  -2.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.68'
$ws.Range("E10").Value = '  -2.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0913'
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.42'
$ws.Range("E12").Value = '  -1.98%  '
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.969'
$ws.Range("E14").Value = '  -4.63%  '
$ws.Range("D15").Value = '2.712.48'
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("E16").Value = '  -2.13%  '
$ws.Range("D17").Value = '2.350.53'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '44.948.33'
$ws.Range("E18").Value = '  +1.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.93'
$ws.Range("E19").Value = '  +7.37%  '
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.17'
$ws.Range("E21").Value = '  -5.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.93'
$ws.Range("E22").Value = '  -2.42%  '
$ws.Range("E23").Value = '  -0.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '258.85'
$ws.Range("E24").Value = '  -2.77%  '
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.99'
$ws.Range("E27").Value = '  -1.82%  '
$ws.Range("E28").Value = '  -6.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.32'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0972'
$ws.Range("E30").Value = '  +5.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.22'
$ws.Range("E31").Value = '  -1.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '36.96'
$ws.Range("E32").Value = '  -6.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '167.43'
$ws.Range("E33").Value = '  -1.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.99'
$ws.Range("E34").Value = '  +4.73%  '
$ws.Range("E35").Value = '  -1.53%  '
$ws.Range("E36").Value = '  -0.37%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.65'
$ws.Range("E37").Value = '  -1.87%  '
$ws.Range("E38").Value = '  +3.65%  '
$ws.Range("E39").Value = '  -3.98%  '
$ws.Range("E40").Value = '  -2.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.77'
$ws.Range("E41").Value = '  +2.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.86'
$ws.Range("E42").Value = '  -3.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.20'
$ws.Range("E43").Value = '  -2.93%  '
$ws.Range("E44").Value = '  -4.46%  '
$ws.Range("E45").Value = '  -0.33%  '
$ws.Range("D46").Value = '1.838.91'
$ws.Range("E46").Value = '  +11.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.74'
$ws.Range("E47").Value = '  -8.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '83.19'
$ws.Range("E48").Value = '  +6.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.66'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.39'
$ws.Range("E50").Value = '  -4.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.13'
$ws.Range("E51").Value = '  +1.13%  '

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
